$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.656.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.702.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9974'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3734'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.93'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.45%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3435'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.181'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07451'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9991'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.87'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.237'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.926'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.705.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001119'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06699'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9981'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '83.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.323'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.95%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.652.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.411'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.758'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.09%  '

$ws.Range("E27").Value = '  +2.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '131.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.894.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.179'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +18.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.743'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.186'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.82%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08797'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.10%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.776'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.94%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '13.66'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.512'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06500'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.918'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02375'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2215'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.274'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6385'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9973'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6076'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.806'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.111'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07265'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.01%  '
